$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.591.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9973"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07035"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.970"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.606"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.655.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001044"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06620"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9967"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.944"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.77%  "

$ws.Range("E22").Value = "  -8.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.585.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.458"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.352"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -15.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.221"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.840.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.764"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08458"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.679"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.41%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.211"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.11%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06022"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02223"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2066"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.120"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.75%  "

$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5896"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.848"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5612"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06952"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.189"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
